# Update the worksheet date header and the 25 division-problem answers
# (two-digit divided by one-digit) throughout the table.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-14 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-15 Thursday", 2) | Out-Null

$d.Content.Find.Execute("97÷6=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "58÷5=11, 3", 2) | Out-Null
$d.Content.Find.Execute("98÷2=49, 0", $true, $false, $false, $false, $false, $true, 1, $false, "61÷5=12, 1", 2) | Out-Null
$d.Content.Find.Execute("31÷7=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "16÷6=2, 4", 2) | Out-Null
$d.Content.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "37÷6=6, 1", 2) | Out-Null
$d.Content.Find.Execute("32÷7=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "96÷4=24, 0", 2) | Out-Null

$d.Content.Find.Execute("65÷9=7, 2", $true, $false, $false, $false, $false, $true, 1, $false, "50÷6=8, 2", 2) | Out-Null
$d.Content.Find.Execute("93÷8=11, 5", $true, $false, $false, $false, $false, $true, 1, $false, "62÷6=10, 2", 2) | Out-Null
$d.Content.Find.Execute("68÷9=7, 5", $true, $false, $false, $false, $false, $true, 1, $false, "51÷6=8, 3", 2) | Out-Null
$d.Content.Find.Execute("13÷3=4, 1", $true, $false, $false, $false, $false, $true, 1, $false, "70÷6=11, 4", 2) | Out-Null
$d.Content.Find.Execute("58÷2=29, 0", $true, $false, $false, $false, $false, $true, 1, $false, "62÷9=6, 8", 2) | Out-Null

$d.Content.Find.Execute("58÷3=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷4=11, 2", 2) | Out-Null
$d.Content.Find.Execute("83÷9=9, 2", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=48, 1", 2) | Out-Null
$d.Content.Find.Execute("69÷9=7, 6", $true, $false, $false, $false, $false, $true, 1, $false, "24÷3=8, 0", 2) | Out-Null
$d.Content.Find.Execute("75÷2=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "78÷5=15, 3", 2) | Out-Null
$d.Content.Find.Execute("47÷3=15, 2", $true, $false, $false, $false, $false, $true, 1, $false, "16÷5=3, 1", 2) | Out-Null

$d.Content.Find.Execute("14÷4=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "46÷7=6, 4", 2) | Out-Null
$d.Content.Find.Execute("90÷6=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "57÷5=11, 2", 2) | Out-Null
$d.Content.Find.Execute("45÷3=15, 0", $true, $false, $false, $false, $false, $true, 1, $false, "25÷5=5, 0", 2) | Out-Null
$d.Content.Find.Execute("60÷9=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "14÷5=2, 4", 2) | Out-Null
$d.Content.Find.Execute("87÷4=21, 3", $true, $false, $false, $false, $false, $true, 1, $false, "46÷4=11, 2", 2) | Out-Null

$d.Content.Find.Execute("82÷6=13, 4", $true, $false, $false, $false, $false, $true, 1, $false, "99÷2=49, 1", 2) | Out-Null
$d.Content.Find.Execute("91÷3=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "58÷6=9, 4", 2) | Out-Null
$d.Content.Find.Execute("46÷9=5, 1", $true, $false, $false, $false, $false, $true, 1, $false, "88÷7=12, 4", 2) | Out-Null
$d.Content.Find.Execute("96÷7=13, 5", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=29, 1", 2) | Out-Null
$d.Content.Find.Execute("14÷8=1, 6", $true, $false, $false, $false, $false, $true, 1, $false, "27÷9=3, 0", 2) | Out-Null
